$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.107.13'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.825.89'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  -0.49%  '
$ws.Range('D5').Value = '''312.57'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').Value = '''0.4581'
$ws.Range('E7').Value = '  +7.60%  '
$ws.Range('E8').Value = '  +1.82%  '
$ws.Range('D9').Value = '''0.07328'
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').Value = '''0.8605'
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('D11').Value = '''21.00'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.826.35'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '''6.694'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').Value = '''93.02'
$ws.Range('E14').Value = '  +6.04%  '
$ws.Range('D15').Value = '''5.345'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '''0.07072'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').Value = '''0.000008833'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '''1.000'
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').Value = '''15.03'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = '27.052.62'
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').Value = '''5.197'
$ws.Range('E22').Value = '  +1.63%  '
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('D24').Value = '''2.002'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('D25').Value = '''151.47'
$ws.Range('E25').Value = '  -0.94%  '
$ws.Range('D26').Value = '''2.228'
$ws.Range('E26').Value = '  +5.69%  '
$ws.Range('D27').Value = '''18.54'
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').Value = '''5.271'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('D29').Value = '''117.54'
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('D30').Value = '''0.08865'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').Value = '''0.7624'
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('D32').Value = '''1.196'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').Value = '''2.966'
$ws.Range('E33').Value = '  +4.98%  '
$ws.Range('D34').Value = '''4.474'
$ws.Range('D35').Value = '''0.9998'
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('D36').Value = '''1.104'
$ws.Range('E36').Value = '  -0.87%  '
$ws.Range('E37').Value = '  +0.67%  '
$ws.Range('D38').Value = '''0.05294'
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').Value = '''0.5353'
$ws.Range('E39').Value = '  +6.87%  '
$ws.Range('D40').Value = '''7.173'
$ws.Range('E40').Value = '  +1.88%  '
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('D42').Value = '''0.1712'
$ws.Range('E42').Value = '  +2.26%  '
$ws.Range('D43').Value = '''0.5219'
$ws.Range('E43').Value = '  +11.38%  '
$ws.Range('D44').Value = '''8.619'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '''10.70'
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('D46').Value = '''1.982'
$ws.Range('E46').Value = '  +10.37%  '
$ws.Range('D47').Value = '''106.05'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.06515'
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''1.677'
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '''0.9997'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '''0.9224'
$ws.Range('E51').Value = '  +1.36%  '
